$d = $word.ActiveDocument

function Split-At($pos, $idx) {
    $rng = $d.Range($pos, $pos)
    $name = "TmpSplit$idx"
    $bm = $d.Bookmarks.Add($name, $rng)
    $d.Bookmarks($name).Delete()
}

# --- Paragraph 2: "In order to directly test the hypothesis ..." ---
$para2 = $d.Paragraphs(2)
$p2start = $para2.Range.Start
$p2end = $para2.Range.End

$fullText = "In order to directly test the hypothesis that higher levels of unemployment directly lead to higher levels of crime, we used ordinary least squares linear regression analysis using base R.  We performed this analysis at the LAD level, and all variables were at the LAD level when this analysis was performed.  Total number of crimes in an LAD was the dependent variable.  The independent variables were the unemployment rate of people between the ages of 16 and 64, the population size of the LAD, and a dummy variable for each year in the study except 2011 (omitted to avoid collinearity)."

$rngAll = $d.Range($p2start, $p2end)
$rngAll.Text = $fullText

# Recompute paragraph 2 start (should be unchanged, but be safe)
$para2 = $d.Paragraphs(2)
$p2start = $para2.Range.Start

$off1 = 125  # after "...we used "
$off2 = 126  # after "o"
$off3 = 153  # after "...linea"
$off4 = 181  # after "...using "
$off5 = 186  # after "base "

Split-At ($p2start + $off1) 1
Split-At ($p2start + $off2) 2
Split-At ($p2start + $off3) 3
Split-At ($p2start + $off4) 4
Split-At ($p2start + $off5) 5

# Insert the permanent bookmark _GoBack at offset 186 (between "base " and "R.")
$bmRng = $d.Range($p2start + $off5, $p2start + $off5)
$d.Bookmarks.Add("_GoBack", $bmRng)

Write-Output "Paragraph 2 updated."
